# Update natmi LR-pair edge statistics (Il17f-Il17ra) with newly recomputed TPM
# values. Only the numeric metric columns E:T change; the Sending
# cluster / Ligand / Receptor / Target cluster label columns (A:D) are
# untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$data = New-Object 'object[,]' 18,16
$data[0,0] = 1
$data[0,1] = 0.5
$data[0,2] = 0.1573005
$data[0,3] = 0.314601
$data[0,4] = 0.2275756804884253
$data[0,5] = 0.1641708874697203
$data[0,6] = 2
$data[0,7] = 1
$data[0,8] = 1.814919
$data[0,9] = 3.629838
$data[0,10] = 0.0191189748917995
$data[0,11] = 0.01299372757260387
$data[0,12] = 0.2854876661595
$data[0,13] = 1.141950664638
$data[0,14] = 0.00435101372124239
$data[0,15] = 0.002133191787134152
$data[1,0] = 1
$data[1,1] = 0.5
$data[1,2] = 0.1573005
$data[1,3] = 0.314601
$data[1,4] = 0.2275756804884253
$data[1,5] = 0.1641708874697203
$data[1,6] = 3
$data[1,7] = 1
$data[1,8] = 12.885218
$data[1,9] = 38.655654
$data[1,10] = 0.1357372750064124
$data[1,11] = 0.1383756071804954
$data[1,12] = 2.026851234009
$data[1,13] = 12.161107404054
$data[1,14] = 0.03089050272722883
$data[1,15] = 0.02271724623498333
$data[2,0] = 1
$data[2,1] = 0.5
$data[2,2] = 0.1573005
$data[2,3] = 0.314601
$data[2,4] = 0.2275756804884253
$data[2,5] = 0.1641708874697203
$data[2,6] = 3
$data[2,7] = 1
$data[2,8] = 20.07571666666667
$data[2,9] = 60.22715
$data[2,10] = 0.2114844369830725
$data[2,11] = 0.2155950705167418
$data[2,12] = 3.157920269525
$data[2,13] = 18.94752161715
$data[2,14] = 0.04812871465913422
$data[2,15] = 0.03539443406083044
$data[3,0] = 1
$data[3,1] = 0.5
$data[3,2] = 0.1573005
$data[3,3] = 0.314601
$data[3,4] = 0.2275756804884253
$data[3,5] = 0.1641708874697203
$data[3,6] = 2
$data[3,7] = 1
$data[3,8] = 3.614881
$data[3,9] = 7.229762000000001
$data[3,10] = 0.03808038765137348
$data[3,11] = 0.02588037202838355
$data[3,12] = 0.5686225887405001
$data[3,13] = 2.274490354962
$data[3,14] = 0.008666170133024348
$data[3,15] = 0.004248803643946254
$data[4,0] = 1
$data[4,1] = 0.5
$data[4,2] = 0.1573005
$data[4,3] = 0.314601
$data[4,4] = 0.2275756804884253
$data[4,5] = 0.1641708874697203
$data[4,6] = 3
$data[4,7] = 1
$data[4,8] = 46.49705
$data[4,9] = 139.49115
$data[4,10] = 0.4898157611952634
$data[4,11] = 0.4993363345386824
$data[4,12] = 7.314009213525001
$data[4,13] = 43.88405528115
$data[4,14] = 0.1114701551679681
$data[4,15] = 0.08197648918709266
$data[5,0] = 1
$data[5,1] = 0.5
$data[5,2] = 0.1573005
$data[5,3] = 0.314601
$data[5,4] = 0.2275756804884253
$data[5,5] = 0.1641708874697203
$data[5,6] = 3
$data[5,7] = 1
$data[5,8] = 10.03984666666667
$data[5,9] = 30.11954
$data[5,10] = 0.1057631642720788
$data[5,11] = 0.107818888163093
$data[5,12] = 1.57927290059
$data[5,13] = 9.47563740354
$data[5,14] = 0.02406912407982744
$data[5,15] = 0.0177007225557335
$data[6,0] = 2
$data[6,1] = 0.6666666666666666
$data[6,2] = 0.4470106666666667
$data[6,3] = 1.341032
$data[6,4] = 0.6467160412856369
$data[6,5] = 0.6998020145050207
$data[6,6] = 2
$data[6,7] = 1
$data[6,8] = 1.814919
$data[6,9] = 3.629838
$data[6,10] = 0.0191189748917995
$data[6,11] = 0.01299372757260387
$data[6,12] = 0.811288152136
$data[6,13] = 4.867728912816
$data[6,14] = 0.01236454775546406
$data[6,15] = 0.009093036731237618
$data[7,0] = 2
$data[7,1] = 0.6666666666666666
$data[7,2] = 0.4470106666666667
$data[7,3] = 1.341032
$data[7,4] = 0.6467160412856369
$data[7,5] = 0.6998020145050207
$data[7,6] = 3
$data[7,7] = 1
$data[7,8] = 12.885218
$data[7,9] = 38.655654
$data[7,10] = 0.1357372750064124
$data[7,11] = 0.1383756071804954
$data[7,12] = 5.759829888325333
$data[7,13] = 51.838468994928
$data[7,14] = 0.08778347314704686
$data[7,15] = 0.09683552866326607
$data[8,0] = 2
$data[8,1] = 0.6666666666666666
$data[8,2] = 0.4470106666666667
$data[8,3] = 1.341032
$data[8,4] = 0.6467160412856369
$data[8,5] = 0.6998020145050207
$data[8,6] = 3
$data[8,7] = 1
$data[8,8] = 20.07571666666667
$data[8,9] = 60.22715
$data[8,10] = 0.2114844369830725
$data[8,11] = 0.2155950705167418
$data[8,12] = 8.974059490977778
$data[8,13] = 80.7665354188
$data[8,14] = 0.1367703778792144
$data[8,15] = 0.1508738646649679
$data[9,0] = 2
$data[9,1] = 0.6666666666666666
$data[9,2] = 0.4470106666666667
$data[9,3] = 1.341032
$data[9,4] = 0.6467160412856369
$data[9,5] = 0.6998020145050207
$data[9,6] = 2
$data[9,7] = 1
$data[9,8] = 3.614881
$data[9,9] = 7.229762000000001
$data[9,10] = 0.03808038765137348
$data[9,11] = 0.02588037202838355
$data[9,12] = 1.615890365730667
$data[9,13] = 9.695342194384
$data[9,14] = 0.0246271975525187
$data[9,15] = 0.0181111364816022
$data[10,0] = 2
$data[10,1] = 0.6666666666666666
$data[10,2] = 0.4470106666666667
$data[10,3] = 1.341032
$data[10,4] = 0.6467160412856369
$data[10,5] = 0.6998020145050207
$data[10,6] = 3
$data[10,7] = 1
$data[10,8] = 46.49705
$data[10,9] = 139.49115
$data[10,10] = 0.4898157611952634
$data[10,11] = 0.4993363345386824
$data[10,12] = 20.78467731853333
$data[10,13] = 187.0620958668
$data[10,14] = 0.3167717100395116
$data[10,15] = 0.3494365728257229
$data[11,0] = 2
$data[11,1] = 0.6666666666666666
$data[11,2] = 0.4470106666666667
$data[11,3] = 1.341032
$data[11,4] = 0.6467160412856369
$data[11,5] = 0.6998020145050207
$data[11,6] = 3
$data[11,7] = 1
$data[11,8] = 10.03984666666667
$data[11,9] = 30.11954
$data[11,10] = 0.1057631642720788
$data[11,11] = 0.107818888163093
$data[11,12] = 4.487918551697778
$data[11,13] = 40.39126696528
$data[11,14] = 0.0683987349118813
$data[11,15] = 0.075451875138224
$data[12,0] = 1
$data[12,1] = 0.3333333333333333
$data[12,2] = 0.08688966666666666
$data[12,3] = 0.260669
$data[12,4] = 0.1257082782259377
$data[12,5] = 0.1360270980252591
$data[12,6] = 2
$data[12,7] = 1
$data[12,8] = 1.814919
$data[12,9] = 3.629838
$data[12,10] = 0.0191189748917995
$data[12,11] = 0.01299372757260387
$data[12,12] = 0.157697706937
$data[12,13] = 0.9461862416219999
$data[12,14] = 0.002403413415093049
$data[12,15] = 0.001767499054232098
$data[13,0] = 1
$data[13,1] = 0.3333333333333333
$data[13,2] = 0.08688966666666666
$data[13,3] = 0.260669
$data[13,4] = 0.1257082782259377
$data[13,5] = 0.1360270980252591
$data[13,6] = 3
$data[13,7] = 1
$data[13,8] = 12.885218
$data[13,9] = 38.655654
$data[13,10] = 0.1357372750064124
$data[13,11] = 0.1383756071804954
$data[13,12] = 1.119592296947333
$data[13,13] = 10.076330672526
$data[13,14] = 0.01706329913213671
$data[13,15] = 0.01882283228224598
$data[14,0] = 1
$data[14,1] = 0.3333333333333333
$data[14,2] = 0.08688966666666666
$data[14,3] = 0.260669
$data[14,4] = 0.1257082782259377
$data[14,5] = 0.1360270980252591
$data[14,6] = 3
$data[14,7] = 1
$data[14,8] = 20.07571666666667
$data[14,9] = 60.22715
$data[14,10] = 0.2114844369830725
$data[14,11] = 0.2155950705167418
$data[14,12] = 1.744372329261111
$data[14,13] = 15.69935096335
$data[14,14] = 0.02658534444472386
$data[14,15] = 0.02932677179094348
$data[15,0] = 1
$data[15,1] = 0.3333333333333333
$data[15,2] = 0.08688966666666666
$data[15,3] = 0.260669
$data[15,4] = 0.1257082782259377
$data[15,5] = 0.1360270980252591
$data[15,6] = 2
$data[15,7] = 1
$data[15,8] = 3.614881
$data[15,9] = 7.229762000000001
$data[15,10] = 0.03808038765137348
$data[15,11] = 0.02588037202838355
$data[15,12] = 0.3140958051296667
$data[15,13] = 1.884574830778
$data[15,14] = 0.004787019965830419
$data[15,15] = 0.003520431902835102
$data[16,0] = 1
$data[16,1] = 0.3333333333333333
$data[16,2] = 0.08688966666666666
$data[16,3] = 0.260669
$data[16,4] = 0.1257082782259377
$data[16,5] = 0.1360270980252591
$data[16,6] = 3
$data[16,7] = 1
$data[16,8] = 46.49705
$data[16,9] = 139.49115
$data[16,10] = 0.4898157611952634
$data[16,11] = 0.4993363345386824
$data[16,12] = 4.040113175483333
$data[16,13] = 36.36101857935
$data[16,14] = 0.06157389598778362
$data[16,15] = 0.0679232725258669
$data[17,0] = 1
$data[17,1] = 0.3333333333333333
$data[17,2] = 0.08688966666666666
$data[17,3] = 0.260669
$data[17,4] = 0.1257082782259377
$data[17,5] = 0.1360270980252591
$data[17,6] = 3
$data[17,7] = 1
$data[17,8] = 10.03984666666667
$data[17,9] = 30.11954
$data[17,10] = 0.1057631642720788
$data[17,11] = 0.107818888163093
$data[17,12] = 0.8723589302511111
$data[17,13] = 7.85123037226
$data[17,14] = 0.01329530528037003
$data[17,15] = 0.01466629046913549
$ws.Range("E2:T19").Value = $data
